$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: numeric value update
$ws.Range("E2").Value = 83.05

# D3 / E3: text values that look numeric -> use leading apostrophe to force text
$ws.Range("D3").Value = "'75.61"
$ws.Range("E3").Value = "'75.61"

# D9 / E9: text values that look numeric -> force text
$ws.Range("D9").Value = "'100.0"
$ws.Range("E9").Value = "'100.0"

# Rows 36, 38, 39, 41-52: D column becomes numeric 100, E column becomes text "Missing"
$rows = @(36, 38, 39, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = 100
    $ws.Range("E$r").Value = "Missing"
}
